$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 12.26167578201206
$ws.Range("D2").Value = 6.92041971848546
$ws.Range("E2").Value = 12.51969760047459
$ws.Range("F2").Value = 43.77443302810957
$ws.Range("G2").Value = 3.706733162025666
$ws.Range("J2").Value = 10.33881531217994
$ws.Range("L2").Value = 8.904533062626106
$ws.Range("M2").Value = 28.90422148313138
$ws.Range("N2").Value = 17.59906630164805
$ws.Range("O2").Value = 34.35877857056661
$ws.Range("C3").Value = 12.2942684492433
$ws.Range("D3").Value = 6.930697108223375
$ws.Range("E3").Value = 12.572098815751
$ws.Range("F3").Value = 43.70184458371285
$ws.Range("G3").Value = 3.710491804684676
$ws.Range("J3").Value = 10.37817805486542
$ws.Range("L3").Value = 8.916655162995491
$ws.Range("M3").Value = 28.29775905681328
$ws.Range("N3").Value = 17.3548014769141
$ws.Range("O3").Value = 34.30373234463958
$ws.Range("C4").Value = 12.31703086550744
$ws.Range("D4").Value = 6.937384166396216
$ws.Range("E4").Value = 12.6058539768777
$ws.Range("F4").Value = 43.66956099064571
$ws.Range("G4").Value = 3.71292084098239
$ws.Range("J4").Value = 10.40351957965715
$ws.Range("L4").Value = 8.924560622737415
$ws.Range("M4").Value = 27.92018390495219
$ws.Range("N4").Value = 17.20521363899108
$ws.Range("O4").Value = 34.2792806680692
$ws.Range("C5").Value = 12.32699586121215
$ws.Range("D5").Value = 6.940204097326682
$ws.Range("E5").Value = 12.62000799232832
$ws.Range("F5").Value = 43.6594968564849
$ws.Range("G5").Value = 3.713941286829352
$ws.Range("J5").Value = 10.41414216532291
$ws.Range("L5").Value = 8.927898770879603
$ws.Range("M5").Value = 27.76519623204183
$ws.Range("N5").Value = 17.14442312401423
$ws.Range("O5").Value = 34.27166865970189
$ws.Range("C6").Value = 12.32869208837576
$ws.Range("D6").Value = 6.940678081675329
$ws.Range("E6").Value = 12.62238235615892
$ws.Range("F6").Value = 43.65801245214163
$ws.Range("G6").Value = 3.714112582209747
$ws.Range("J6").Value = 10.4159239223041
$ws.Range("L6").Value = 8.928460120270813
$ws.Range("M6").Value = 27.73939841548325
$ws.Range("N6").Value = 17.13434106690744
$ws.Range("O6").Value = 34.27054677285725
$ws.Range("C7").Value = 12.3171624700322
$ws.Range("D7").Value = 6.937421812479467
$ws.Range("E7").Value = 12.60604324792844
$ws.Range("F7").Value = 43.66941274354091
$ws.Range("G7").Value = 3.712934479049267
$ws.Range("J7").Value = 10.40366164115617
$ws.Range("L7").Value = 8.924605169603772
$ws.Range("M7").Value = 27.91809798185267
$ws.Range("N7").Value = 17.20439302593033
$ws.Range("O7").Value = 34.27916848504053
$ws.Range("C8").Value = 12.27234142555849
$ws.Range("D8").Value = 6.923885284769137
$ws.Range("E8").Value = 12.53743812476103
$ws.Range("F8").Value = 43.74685683412237
$ws.Range("G8").Value = 3.708004052362274
$ws.Range("J8").Value = 10.35214471084487
$ws.Range("L8").Value = 8.908616978670873
$ws.Range("M8").Value = 28.69630278493509
$ws.Range("N8").Value = 17.51480387154522
$ws.Range("O8").Value = 34.3378604656446
$ws.Range("C9").Value = 12.20637823547111
$ws.Range("D9").Value = 6.900321090366153
$ws.Range("E9").Value = 12.41539604154928
$ws.Range("F9").Value = 43.99597661949608
$ws.Range("G9").Value = 3.699292113186764
$ws.Range("J9").Value = 10.26038573328713
$ws.Range("L9").Value = 8.880918778862746
$ws.Range("M9").Value = 30.17328737062231
$ws.Range("N9").Value = 18.12361779106923
$ws.Range("O9").Value = 34.52694010779736
$ws.Range("C10").Value = 12.1714269117747
$ws.Range("D10").Value = 6.884814458536052
$ws.Range("E10").Value = 12.33327772017182
$ws.Range("F10").Value = 44.23779054514355
$ws.Range("G10").Value = 3.693467399954912
$ws.Range("J10").Value = 10.19856381361709
$ws.Range("L10").Value = 8.8627765167207
$ws.Range("M10").Value = 31.21871709625103
$ws.Range("N10").Value = 18.56701910888594
$ws.Range("O10").Value = 34.71058109564893
$ws.Range("C11").Value = 12.15848892624992
$ws.Range("D11").Value = 6.8781498764556
$ws.Range("E11").Value = 12.29754367561165
$ws.Range("F11").Value = 44.36039167568377
$ws.Range("G11").Value = 3.690941106003895
$ws.Range("J11").Value = 10.17164221984987
$ws.Range("L11").Value = 8.854998170078082
$ws.Range("M11").Value = 31.68380379016588
$ws.Range("N11").Value = 18.76708066837963
$ws.Range("O11").Value = 34.80370778036674
$ws.Range("C12").Value = 12.15401763343427
$ws.Range("D12").Value = 6.87568200427477
$ws.Range("E12").Value = 12.28424427329026
$ws.Range("F12").Value = 44.40860975339569
$ws.Range("G12").Value = 3.690002090090879
$ws.Range("J12").Value = 10.16161961920428
$ws.Range("L12").Value = 8.852120634887093
$ws.Range("M12").Value = 31.85826867606979
$ws.Range("N12").Value = 18.84253765501906
$ws.Range("O12").Value = 34.84033658690053
$ws.Range("C13").Value = 12.15496153586191
$ws.Range("D13").Value = 6.876211022561937
$ws.Range("E13").Value = 12.28709822049741
$ws.Range("F13").Value = 44.39814579157238
$ws.Range("G13").Value = 3.690203541276246
$ws.Range("J13").Value = 10.16377052672917
$ws.Range("L13").Value = 8.852737345972052
$ws.Range("M13").Value = 31.82077028506093
$ws.Range("N13").Value = 18.82630109180154
$ws.Range("O13").Value = 34.83238752547308
$ws.Range("C14").Value = 12.1581124797222
$ws.Range("D14").Value = 6.877945724414217
$ws.Range("E14").Value = 12.29644487536601
$ws.Range("F14").Value = 44.36432282290229
$ws.Range("G14").Value = 3.690863499747191
$ws.Range("J14").Value = 10.17081421115895
$ws.Range("L14").Value = 8.854760073300064
$ws.Range("M14").Value = 31.69819105493947
$ws.Range("N14").Value = 18.77329497924772
$ws.Range("O14").Value = 34.80669401591241
$ws.Range("C15").Value = 12.16009832975272
$ws.Range("D15").Value = 6.879015549840449
$ws.Range("E15").Value = 12.30220019670976
$ws.Range("F15").Value = 44.34383796475911
$ws.Range("G15").Value = 3.69127003700981
$ws.Range("J15").Value = 10.17515105146487
$ws.Range("L15").Value = 8.856007893188563
$ws.Range("M15").Value = 31.62288827338979
$ws.Range("N15").Value = 18.74078597272719
$ws.Range("O15").Value = 34.79113312015333
$ws.Range("C16").Value = 12.17233222233489
$ws.Range("D16").Value = 6.885257829274792
$ws.Range("E16").Value = 12.3356455824181
$ws.Range("F16").Value = 44.23003015105655
$ws.Range("G16").Value = 3.693634973065602
$ws.Range("J16").Value = 10.2003473159748
$ws.Range("L16").Value = 8.863294375682278
$ws.Range("M16").Value = 31.18809892316285
$ws.Range("N16").Value = 18.55390587355626
$ws.Range("O16").Value = 34.70468674755413
$ws.Range("C17").Value = 12.18059743093215
$ws.Range("D17").Value = 6.88918691302563
$ws.Range("E17").Value = 12.35657800510038
$ws.Range("F17").Value = 44.16342631269034
$ws.Range("G17").Value = 3.695117313025039
$ws.Range("J17").Value = 10.21611160511442
$ws.Range("L17").Value = 8.867885754671697
$ws.Range("M17").Value = 30.91857748413999
$ws.Range("N17").Value = 18.43879251187378
$ws.Range("O17").Value = 34.65410036941955
$ws.Range("C18").Value = 12.18562999324608
$ws.Range("D18").Value = 6.891483484572265
$ws.Range("E18").Value = 12.36877050583476
$ws.Range("F18").Value = 44.12630479168168
$ws.Range("G18").Value = 3.695981536907263
$ws.Range("J18").Value = 10.22529194342613
$ws.Range("L18").Value = 8.870571289441608
$ws.Range("M18").Value = 30.76257708796989
$ws.Range("N18").Value = 18.37243139514892
$ws.Range("O18").Value = 34.62590787649101
$ws.Range("C19").Value = 12.18738171942732
$ws.Range("D19").Value = 6.892267366169111
$ws.Range("E19").Value = 12.3729249363883
$ws.Range("F19").Value = 44.11394055311846
$ws.Range("G19").Value = 3.696276147492451
$ws.Range("J19").Value = 10.22841970251824
$ws.Range("L19").Value = 8.871488250929064
$ws.Range("M19").Value = 30.70959453292465
$ws.Range("N19").Value = 18.34993877232939
$ws.Range("O19").Value = 34.61651797009494
$ws.Range("C20").Value = 12.17968873178835
$ws.Range("D20").Value = 6.888764861553937
$ws.Range("E20").Value = 12.3543339133941
$ws.Range("F20").Value = 44.17039366542118
$ws.Range("G20").Value = 3.694958313447073
$ws.Range("J20").Value = 10.21442176598612
$ws.Range("L20").Value = 8.867392370970832
$ws.Range("M20").Value = 30.94737089084858
$ws.Range("N20").Value = 18.45106263193497
$ws.Range("O20").Value = 34.65939197309606
$ws.Range("C21").Value = 12.15717533752953
$ws.Range("D21").Value = 6.877434685527987
$ws.Range("E21").Value = 12.2936932389809
$ws.Range("F21").Value = 44.37420899677127
$ws.Range("G21").Value = 3.690669176144292
$ws.Range("J21").Value = 10.1687406476162
$ws.Range("L21").Value = 8.854164107328625
$ws.Range("M21").Value = 31.73424152780891
$ws.Range("N21").Value = 18.78887286568135
$ws.Range("O21").Value = 34.81420394034361
$ws.Range("C22").Value = 12.14495740102105
$ws.Range("D22").Value = 6.870355292949995
$ws.Range("E22").Value = 12.2554146765045
$ws.Range("F22").Value = 44.51784609461532
$ws.Range("G22").Value = 3.687968728022623
$ws.Range("J22").Value = 10.13988773161819
$ws.Range("L22").Value = 8.845914633257493
$ws.Range("M22").Value = 32.23880963547492
$ws.Range("N22").Value = 19.00786254089702
$ws.Range("O22").Value = 34.92332305638427
$ws.Range("C23").Value = 12.15124928742453
$ws.Range("D23").Value = 6.874103958574318
$ws.Range("E23").Value = 12.27572111534361
$ws.Range("F23").Value = 44.4402371108501
$ws.Range("G23").Value = 3.689400641703994
$ws.Range("J23").Value = 10.15519560886021
$ws.Range("L23").Value = 8.850281400936039
$ws.Range("M23").Value = 31.97044492750812
$ws.Range("N23").Value = 18.89116833288719
$ws.Range("O23").Value = 34.8643630584578
$ws.Range("C24").Value = 12.18009868024032
$ws.Range("D24").Value = 6.888955553624579
$ws.Range("E24").Value = 12.3553479744865
$ws.Range("F24").Value = 44.16724007953633
$ws.Range("G24").Value = 3.695030159747187
$ws.Range("J24").Value = 10.21518537700309
$ws.Range("L24").Value = 8.867615286772839
$ws.Range("M24").Value = 30.93435664763202
$ws.Range("N24").Value = 18.44551586913838
$ws.Range("O24").Value = 34.65699686474304
$ws.Range("C25").Value = 12.22186006944816
$ws.Range("D25").Value = 6.906377876489529
$ws.Range("E25").Value = 12.44708150619202
$ws.Range("F25").Value = 43.9182110909229
$ws.Range("G25").Value = 3.701547260687791
$ws.Range("J25").Value = 10.28422275110018
$ws.Range("L25").Value = 8.888022734164121
$ws.Range("M25").Value = 29.78002504549185
$ws.Range("N25").Value = 17.95932074830214
$ws.Range("O25").Value = 34.46789852867372
